$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4475
$ws.Range("K3").Value = 4591
$ws.Range("K4").Value = 925
$ws.Range("K5").Value = 332
$ws.Range("K6").Value = 5182
$ws.Range("K7").Value = 15505

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 95
$ws.Range("K7").Value = 206

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 289
$ws.Range("K3").Value = 311
$ws.Range("K6").Value = 350
$ws.Range("K7").Value = 1036

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 115
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 180
$ws.Range("K3").Value = 243
$ws.Range("K7").Value = 654

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 146
$ws.Range("K5").Value = 26
$ws.Range("K7").Value = 525

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 108
$ws.Range("K3").Value = 88
$ws.Range("K6").Value = 139
$ws.Range("K7").Value = 355

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 108
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 134
$ws.Range("K4").Value = 60
$ws.Range("K7").Value = 455
$ws.Range("K8").Value = 1036
$ws.Range("K18").Value = 105
$ws.Range("K19").Value = 467
$ws.Range("K20").Value = 357
$ws.Range("K23").Value = 158
$ws.Range("K24").Value = 45
$ws.Range("K25").Value = 75
$ws.Range("K27").Value = 142
$ws.Range("K29").Value = 820
$ws.Range("K31").Value = 170
$ws.Range("K33").Value = 654
$ws.Range("K34").Value = 81
$ws.Range("K36").Value = 194
$ws.Range("K37").Value = 525
$ws.Range("K42").Value = 573
$ws.Range("K43").Value = 139
$ws.Range("K46").Value = 34
$ws.Range("K50").Value = 83
$ws.Range("K52").Value = 408
$ws.Range("K53").Value = 206
$ws.Range("K54").Value = 292
$ws.Range("K60").Value = 101
$ws.Range("J63").Value = 111
$ws.Range("K63").Value = 52
$ws.Range("K65").Value = 355
$ws.Range("K67").Value = 596
$ws.Range("K68").Value = 40
$ws.Range("K71").Value = 50
$ws.Range("K72").Value = 71
$ws.Range("K76").Value = 215
$ws.Range("K77").Value = 112
$ws.Range("J78").Value = 340
$ws.Range("K79").Value = 388
$ws.Range("K83").Value = 331
$ws.Range("K85").Value = 694
$ws.Range("K88").Value = 181
$ws.Range("K89").Value = 223
$ws.Range("K91").Value = 168
$ws.Range("K92").Value = 59
$ws.Range("K94").Value = 192
$ws.Range("K96").Value = 173
$ws.Range("K99").Value = 262
$ws.Range("K101").Value = 15505

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 212
$ws.Range("K7").Value = 596

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 80
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 292

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 236
$ws.Range("K3").Value = 291
$ws.Range("K7").Value = 820

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 143
$ws.Range("K7").Value = 467

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 44
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 154
$ws.Range("K3").Value = 177
$ws.Range("K7").Value = 573

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J5").Value = 6
$ws.Range("J7").Value = 340

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 42
$ws.Range("K3").Value = 80
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 128
$ws.Range("K7").Value = 388

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 113
$ws.Range("K7").Value = 357

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 76
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 455

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 56
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 192

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 43
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 70
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 223

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 67

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 242
$ws.Range("K6").Value = 161
$ws.Range("K7").Value = 694

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 43
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K6").Value = 154
$ws.Range("K7").Value = 408

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 60
